$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.484.12"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "'2.035.79"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'229.73"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'56.24"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "'0.0804"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "'2.340.01"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "'14.42"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "'5.22"
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "'2.033.71"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'37.379.62"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "'0.0₃0827"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'223.65"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").Value = "'2.26"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").Value = "'164.69"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'0.132"
$ws.Range("E28").Value = "  +6.80%  "
$ws.Range("D29").Value = "'18.79"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").Value = "'4.49"
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").Value = "'2.01"
$ws.Range("E35").Value = "  +9.46%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +9.37%  "
$ws.Range("D38").Value = "'3.23"
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'1.476.30"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "'94.89"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.23"
$ws.Range("E45").Value = "  +17.35%  "
$ws.Range("D46").Value = "'16.30"
$ws.Range("E46").Value = "  -4.50%  "
$ws.Range("D47").Value = "'1.11"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.95"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'7.11"
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").Value = "'2.228.00"
$ws.Range("E51").Value = "  +0.91%  "
